$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new column before column A - this shifts existing data (and column widths) one
# column to the right (A->B, B->C, ... X->Y).
$ws.Columns.Item(1).Insert()

# New column A holds the laboratory number for each sample; give it a header and a narrow width.
$ws.Range("A1").Value2 = "Lab. #"
$ws.Columns.Item(1).ColumnWidth = 6.8

# Fill in the laboratory numbers for each data row.
$labNumbers = @{
    2  = 10815
    3  = 10989
    4  = 10815
    5  = 10990
    6  = 10815
    7  = 10991
    8  = 10815
    9  = 10992
    10 = 10815
    11 = 10993
    12 = 10815
}

foreach ($row in $labNumbers.Keys) {
    $ws.Cells.Item($row, 1).Value2 = $labNumbers[$row]
}

# Rows that belong to the same laboratory sample (10815) are highlighted with a light green fill
# across the whole row of data.
$fillRows = @(2, 4, 6, 8, 10, 12)
foreach ($row in $fillRows) {
    $rng = $ws.Range($ws.Cells.Item($row, 1), $ws.Cells.Item($row, 25))
    $rng.Interior.Color = 12379352
}
